# Applies the Mon Jun 19 05:39:44 UTC 2023 cryptos-list refresh described by the diff.
# Column D ("Price") holds numeric-looking text (e.g. "138.00", "1.230") that must stay
# text (trailing zeros / thousand-dot formatting matter), so those assignments are
# prefixed with a leading apostrophe -- the same "force text" trick Excel itself uses
# when a typed value looks like a number but the cell should remain text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.427.95'
$ws.Range('E2').Value = '  -0.60%  '

$ws.Range('D3').Value = '1.724.60'
$ws.Range('E3').Value = '  -0.42%  '

$ws.Range('D4').Value = '''0.9995'
$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = '''243.08'
$ws.Range('E5').Value = '  -0.89%  '

$ws.Range('D6').Value = '''0.9998'
$ws.Range('E6').Value = '  +0.08%  '

$ws.Range('D7').Value = '''0.4911'
$ws.Range('E7').Value = '  +2.05%  '

$ws.Range('E8').Value = '  -2.36%  '

$ws.Range('D9').Value = '''0.06202'
$ws.Range('E9').Value = '  +0.11%  '

$ws.Range('D10').Value = '1.724.02'
$ws.Range('E10').Value = '  -0.39%  '

$ws.Range('D11').Value = '''0.06996'
$ws.Range('E11').Value = '  -2.75%  '

$ws.Range('D12').Value = '''15.45'
$ws.Range('E12').Value = '  -0.89%  '

$ws.Range('D13').Value = '''4.561'
$ws.Range('E13').Value = '  +0.49%  '

$ws.Range('D14').Value = '''0.5990'
$ws.Range('E14').Value = '  -2.08%  '

$ws.Range('D15').Value = '''77.29'

$ws.Range('D16').Value = '''0.9998'
$ws.Range('E16').Value = '  +0.06%  '

$ws.Range('D17').Value = '26.423.55'
$ws.Range('E17').Value = '  -0.58%  '

$ws.Range('D18').Value = '''0.9998'
$ws.Range('E18').Value = '  +0.05%  '

$ws.Range('D19').Value = '''0.000007165'
$ws.Range('E19').Value = '  +2.68%  '

$ws.Range('D20').Value = '''11.35'
$ws.Range('E20').Value = '  -1.89%  '

$ws.Range('D21').Value = '1.944.53'
$ws.Range('E21').Value = '  -0.36%  '

$ws.Range('D22').Value = '''4.479'
$ws.Range('E22').Value = '  -1.17%  '

$ws.Range('D23').Value = '''8.591'
$ws.Range('E23').Value = '  -2.56%  '

$ws.Range('D24').Value = '''5.165'
$ws.Range('E24').Value = '  -1.90%  '

$ws.Range('D25').Value = '''138.00'
$ws.Range('E25').Value = '  +0.67%  '

$ws.Range('D26').Value = '''15.24'
$ws.Range('E26').Value = '  -0.78%  '

$ws.Range('D27').Value = '''1.399'
$ws.Range('E27').Value = '  -0.47%  '

$ws.Range('E28').Value = '  -0.29%  '

$ws.Range('D29').Value = '''1.714'
$ws.Range('E29').Value = '  -4.04%  '

$ws.Range('D30').Value = '''3.948'
$ws.Range('E30').Value = '  -1.08%  '

$ws.Range('D31').Value = '''0.07963'
$ws.Range('E31').Value = '  -0.95%  '

$ws.Range('D32').Value = '''3.675'
$ws.Range('E32').Value = '  -0.74%  '

$ws.Range('D33').Value = '''0.04532'
$ws.Range('E33').Value = '  +0.02%  '

$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').Value = '''0.9991'
$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '''2.601'
$ws.Range('E35').Value = '  -0.57%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '''0.9945'
$ws.Range('E36').Value = '  -1.64%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '''0.6254'
$ws.Range('E37').Value = '  -0.32%  '

$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '''0.9263'
$ws.Range('E38').Value = '  +2.10%  '

$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''1.961'
$ws.Range('E39').Value = '  -5.80%  '

$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '''2.391'
$ws.Range('E40').Value = '  -0.72%  '

$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '''0.9998'
$ws.Range('E41').Value = '  -0.12%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.01486'
$ws.Range('E42').Value = '  -1.35%  '

$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '''99.28'
$ws.Range('E43').Value = '  -3.12%  '

$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''5.340'
$ws.Range('E44').Value = '  -3.39%  '

$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '''0.3847'
$ws.Range('E45').Value = '  -1.28%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '''6.758'
$ws.Range('E46').Value = '  -3.99%  '

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '''0.1167'
$ws.Range('E47').Value = '  -1.39%  '

$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.05367'
$ws.Range('E48').Value = '  -0.27%  '

$ws.Range('D49').Value = '''30.13'
$ws.Range('E49').Value = '  -2.07%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.694'
$ws.Range('E50').Value = '  -2.47%  '

$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '''1.230'
$ws.Range('E51').Value = '  -1.73%  '
